$d = $word.ActiveDocument

# --- Change 1: merge ", " + "if" + " " runs into a single ", if " run -----
# (collapses the run split around the spell-checked word "if" in
# "Similarly, if our model has too many ...")
$r1 = $d.Content
$found1 = $r1.Find.Execute(", if ", $true, $false, $false, $false, $false, $true, 1, $false, ", if ", 2)

# --- Change 2: append two new paragraphs after the "... daily questions?"
# paragraph: one empty paragraph, then a new Q&A paragraph.
$r2 = $d.Content
$found2 = $r2.Find.Execute("daily questions?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.Move(4, 1) | Out-Null
$r2.InsertParagraphAfter()
$r2.Move(4, 1) | Out-Null
$r2.InsertAfter("Q: What is ROC curve?")

Write-Output "change1 found: $found1; change2 found: $found2"
